$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New week column (25_02_2024) header and values
$ws.Range("D1").Value = "25_02_2024"
$ws.Range("D2").Value = 5
$ws.Range("D3").Value = 5
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 1.2
$ws.Range("D6").Value = 4

# Update the selected cell to reflect where the user ended up editing
$ws.Range("D7").Select()
